$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.551.53"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "'1.622.53"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'211.65"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'23.15"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'1.851.85"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "'1.620.17"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'0.549"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'65.27"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'27.498.72"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "'232.29"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'0.0₃0720"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'7.57"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'10.22"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("E24").Value = "  +6.25%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'6.88"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.56"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'3.28"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'1.478.67"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'0.953"
$ws.Range("E37").Value = "  +7.89%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.872"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.554"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.02"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'67.64"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").Value = "'1.762.54"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  -6.47%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "'87.25"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0106"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").Value = "'7.70"
$ws.Range("E51").Value = "  -1.06%  "
